# Filter - Study - Test Suit
# Rename the "CasesTab" row label to "ParticipantsTab" and move the
# active selection to A2 on the "startup" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the tab-name cell in row 2 from "CasesTab" to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Move/refresh the active selection to A2
$ws.Activate()
$ws.Range("A2").Select()
